# This sheet is a weekly price log for "Brócoli" at "Vega Monumental Concepción".
# A new week's pair of rows (Primera / Segunda quality) is inserted in the middle
# of the existing data block (before the row that used to be row 558), which
# shifts every following row down by two and grows the used range from
# A1:R601 to A1:R603.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 558-559, pushing the old 558..601 block down to 560..603.
$ws.Range("A558:A559").EntireRow.Insert()

# Row 558: Brócoli, Primera, week of 45223 (2023-10-24)
$ws.Range("A558").Value2 = 11
$ws.Range("B558").Value2 = 'Vega Monumental Concepción'
$ws.Range("C558").Value2 = 'Bíobío'
$ws.Range("D558").Value2 = 45223
$ws.Range("E558").Value2 = 8
$ws.Range("F558").Value2 = 100112023
$ws.Range("G558").Value2 = 'Brócoli'
$ws.Range("H558").Value2 = 'Sin especificar'
$ws.Range("I558").Value2 = 'Primera'
$ws.Range("J558").Value2 = 1500
$ws.Range("K558").Value2 = 800
$ws.Range("L558").Value2 = 900
$ws.Range("M558").Value2 = 867
$ws.Range("N558").Value2 = '$/unidad'
$ws.Range("O558").Value2 = 'Región Metropolitana'
$ws.Range("P558").Value2 = 867
$ws.Range("Q558").Value2 = 1
$ws.Range("R558").Value2 = 'Hortaliza'

# Row 559: Brócoli, Segunda, week of 45223 (2023-10-24)
$ws.Range("A559").Value2 = 11
$ws.Range("B559").Value2 = 'Vega Monumental Concepción'
$ws.Range("C559").Value2 = 'Bíobío'
$ws.Range("D559").Value2 = 45223
$ws.Range("E559").Value2 = 8
$ws.Range("F559").Value2 = 100112023
$ws.Range("G559").Value2 = 'Brócoli'
$ws.Range("H559").Value2 = 'Sin especificar'
$ws.Range("I559").Value2 = 'Segunda'
$ws.Range("J559").Value2 = 1000
$ws.Range("K559").Value2 = 700
$ws.Range("L559").Value2 = 700
$ws.Range("M559").Value2 = 700
$ws.Range("N559").Value2 = '$/unidad'
$ws.Range("O559").Value2 = 'Región Metropolitana'
$ws.Range("P559").Value2 = 700
$ws.Range("Q559").Value2 = 1
$ws.Range("R559").Value2 = 'Hortaliza'

# Make sure the two new date cells carry the same date number format (style)
# as the rest of column D (style index 2 in this workbook).
$ws.Range("D558").NumberFormat = $ws.Range("D560").NumberFormat
$ws.Range("D559").NumberFormat = $ws.Range("D560").NumberFormat
